# Auto update Excel log 2026-02-04 14:31:29
# Appends new sensor-log rows to the PIR, Humidity and Temperature sheets.

$wb = $excel.ActiveWorkbook

function Add-LogRows {
    param($ws, $rows)

    foreach ($row in $rows) {
        $r = [int]$row[0]
        # Force every cell to text first -- the source data (dates like
        # "2026-02-04" and percentages like "79.5%") would otherwise be
        # auto-coerced into numeric/date values by the COM layer, whereas
        # the log keeps them as literal strings.
        for ($col = 1; $col -le 6; $col++) {
            $ws.Cells.Item($r, $col).NumberFormat = "@"
        }
        $ws.Cells.Item($r, 1).Value = $row[1]
        $ws.Cells.Item($r, 2).Value = $row[2]
        $ws.Cells.Item($r, 3).Value = $row[3]
        $ws.Cells.Item($r, 4).Value = $row[4]
        $ws.Cells.Item($r, 5).Value = $row[5]
        $ws.Cells.Item($r, 6).Value = $row[6]
    }
}

# --- PIR sheet: append rows 350-365 ---
$wsPIR = $wb.Worksheets.Item("PIR")
$pirRows = @(
  @("350","2026-02-04","14:30:25","14:00","Bathroom","No Motion","Inactive"),
  @("351","2026-02-04","14:30:26","14:00","Bathroom","No Motion","Inactive"),
  @("352","2026-02-04","14:30:27","14:00","Bathroom","Motion Detected","Active"),
  @("353","2026-02-04","14:30:30","14:00","Bathroom","No Motion","Inactive"),
  @("354","2026-02-04","14:30:30","14:00","Bathroom","Motion Detected","Active"),
  @("355","2026-02-04","14:30:38","14:00","Bathroom","No Motion","Inactive"),
  @("356","2026-02-04","14:30:43","14:00","Bathroom","No Motion","Inactive"),
  @("357","2026-02-04","14:30:48","14:00","Bathroom","No Motion","Inactive"),
  @("358","2026-02-04","14:30:54","14:00","Bathroom","No Motion","Inactive"),
  @("359","2026-02-04","14:30:58","14:00","Bathroom","No Motion","Inactive"),
  @("360","2026-02-04","14:30:59","14:00","Bathroom","Motion Detected","Active"),
  @("361","2026-02-04","14:31:09","14:00","Bathroom","No Motion","Inactive"),
  @("362","2026-02-04","14:31:14","14:00","Bathroom","No Motion","Inactive"),
  @("363","2026-02-04","14:31:19","14:00","Bathroom","No Motion","Inactive"),
  @("364","2026-02-04","14:31:22","14:00","Bathroom","No Motion","Inactive"),
  @("365","2026-02-04","14:31:23","14:00","Bathroom","Motion Detected","Active")
)
Add-LogRows $wsPIR $pirRows

# --- Humidity sheet: append rows 288-295 ---
$wsHumidity = $wb.Worksheets.Item("Humidity")
$humidityRows = @(
  @("288","2026-02-04","14:30:23","14:00","Bathroom","79.5%","Active"),
  @("289","2026-02-04","14:30:28","14:00","Bathroom","78.6%","Active"),
  @("290","2026-02-04","14:30:32","14:00","Bathroom","78.6%","Active"),
  @("291","2026-02-04","14:30:52","14:00","Bathroom","79.2%","Active"),
  @("292","2026-02-04","14:31:02","14:00","Bathroom","79.5%","Active"),
  @("293","2026-02-04","14:31:07","14:00","Bathroom","79.7%","Active"),
  @("294","2026-02-04","14:31:12","14:00","Bathroom","78.9%","Active"),
  @("295","2026-02-04","14:31:17","14:00","Bathroom","79.9%","Active")
)
Add-LogRows $wsHumidity $humidityRows

# --- Temperature sheet: append rows 288-295 ---
$wsTemperature = $wb.Worksheets.Item("Temperature")
$temperatureRows = @(
  @("288","2026-02-04","14:30:24","14:00","Bathroom","24.3C","Active"),
  @("289","2026-02-04","14:30:29","14:00","Bathroom","24.3C","Active"),
  @("290","2026-02-04","14:30:33","14:00","Bathroom","24.4C","Active"),
  @("291","2026-02-04","14:30:53","14:00","Bathroom","24.3C","Active"),
  @("292","2026-02-04","14:31:03","14:00","Bathroom","24.3C","Active"),
  @("293","2026-02-04","14:31:08","14:00","Bathroom","24.3C","Active"),
  @("294","2026-02-04","14:31:13","14:00","Bathroom","24.3C","Active"),
  @("295","2026-02-04","14:31:18","14:00","Bathroom","24.3C","Active")
)
Add-LogRows $wsTemperature $temperatureRows
